$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# This sheet is a weekly price log where each week a new observation is
# inserted at the top of the data block (row 18) and every later row shifts
# down by one, with the oldest row (198) falling off the bottom onto a new
# last row (199).
#
# Columns D (Fecha), J (Volumen), K (Precio minimo), M (Precio promedio
# ponderado) and P (Precio $/Kg) are the values that move with each record;
# the other columns (A,B,C,E,F,G,H,I,L,N,O,Q,R) are constant across the whole
# table for this subset, so they do not need to be touched.
# ---------------------------------------------------------------------------

$firstRow = 18
$lastRow  = 198
$newRow   = 199

# Snapshot the original values for D, J, K, M, P for every row in the moving
# range (18..198) BEFORE any writes happen, so the shift computation always
# uses the pre-edit data.
$D = @{}
$J = @{}
$K = @{}
$M = @{}
$P = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $D[$r] = $ws.Range("D$r").Value2()
    $J[$r] = $ws.Range("J$r").Value2()
    $K[$r] = $ws.Range("K$r").Value2()
    $M[$r] = $ws.Range("M$r").Value2()
    $P[$r] = $ws.Range("P$r").Value2()
}

# Append the new row (199) as a full copy of the row that is about to fall
# off the bottom of the table (the original row 198), preserving formatting.
$fullRowVals = $ws.Range("A$lastRow`:R$lastRow").Value2()
$ws.Range("A$newRow`:R$newRow").Value = $fullRowVals
$ws.Range("D$newRow").NumberFormat = $ws.Range("D$lastRow").NumberFormat()

# Shift rows 19..198 down from the original row immediately above them
# (row N takes the pre-edit values that used to live in row N-1).
for ($r = ($firstRow + 1); $r -le $lastRow; $r++) {
    $src = $r - 1
    $ws.Range("D$r").Value = $D[$src]
    $ws.Range("J$r").Value = $J[$src]
    $ws.Range("K$r").Value = $K[$src]
    $ws.Range("M$r").Value = $M[$src]
    $ws.Range("P$r").Value = $P[$src]
}

# Row 18 becomes the brand new weekly observation.
$ws.Range("D18").Value = 44552
$ws.Range("J18").Value = 3000
